# Update the division problems in the single table of this worksheet.
# Each data row of the table (rows 1, 5, 9, 13, 17 — separated by blank
# spacer rows) holds 5 "dividend÷divisor=" cells that are replaced with
# new values, per the target diff. We address cells directly by
# (row, column) through the Tables/Cell object model rather than relying
# on document-wide Find/Replace, since several new values collide with
# other cells' old values (e.g. "75÷3=" and "58÷3=") and a naive
# sequential text replace could clobber the wrong cell.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "25÷6="   # was 50÷3=
$t.Cell(1,2).Range.Text  = "75÷3="   # was 20÷7=
$t.Cell(1,3).Range.Text  = "51÷8="   # was 13÷2=
$t.Cell(1,4).Range.Text  = "58÷3="   # was 68÷3=
$t.Cell(1,5).Range.Text  = "67÷9="   # was 40÷8=

$t.Cell(5,1).Range.Text  = "65÷8="   # was 95÷5=
$t.Cell(5,2).Range.Text  = "60÷7="   # was 90÷8=
$t.Cell(5,3).Range.Text  = "70÷2="   # was 88÷7=
$t.Cell(5,4).Range.Text  = "18÷5="   # was 37÷5=
$t.Cell(5,5).Range.Text  = "59÷4="   # was 22÷3=

$t.Cell(9,1).Range.Text  = "75÷2="   # was 52÷9=
$t.Cell(9,2).Range.Text  = "93÷4="   # was 12÷7=
$t.Cell(9,3).Range.Text  = "34÷8="   # was 77÷9=
$t.Cell(9,4).Range.Text  = "42÷9="   # was 30÷4=
$t.Cell(9,5).Range.Text  = "23÷3="   # was 75÷3=

$t.Cell(13,1).Range.Text = "53÷3="   # was 25÷4=
$t.Cell(13,2).Range.Text = "27÷5="   # was 22÷2=
$t.Cell(13,3).Range.Text = "91÷2="   # was 98÷7=
$t.Cell(13,4).Range.Text = "38÷9="   # was 99÷6=
$t.Cell(13,5).Range.Text = "56÷3="   # was 42÷3=

$t.Cell(17,1).Range.Text = "63÷4="   # was 25÷5=
$t.Cell(17,2).Range.Text = "25÷8="   # was 99÷3=
$t.Cell(17,3).Range.Text = "30÷2="   # was 96÷5=
$t.Cell(17,4).Range.Text = "79÷3="   # was 58÷3=
$t.Cell(17,5).Range.Text = "18÷5="   # was 59÷9=
